$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fixed loan balance error: update Repayment Amount (column C) for rows 2-13
$newAmount = 18333.33333333333

for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 3).Value = $newAmount
}
